$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.989.79'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.626.99'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.69'
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.251'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("E9").Value = '  -2.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.44'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -5.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = '1.854.28'
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").Value = '1.622.14'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("E15").Value = '  -2.82%  '
$ws.Range("D16").Value = '25.992.66'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.42'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.85%  '
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '191.91'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.02%  '
$ws.Range("E21").Value = '  -2.01%  '
$ws.Range("E22").Value = '  -3.49%  '
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  +1.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.20'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.41%  '
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.71'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.23'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("E32").Value = '  -3.42%  '
$ws.Range("E33").Value = '  -4.92%  '
$ws.Range("E34").Value = '  -1.14%  '
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").Value = '1.122.07'
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.849'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -5.69%  '
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("E39").Value = '  -3.57%  '
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.35'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '1.764.03'
$ws.Range("E42").Value = '  -0.60%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.753'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.09'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -5.93%  '
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("E46").Value = '  -3.43%  '
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("E51").Value = '  -3.14%  '
